$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 11, shifting existing rows 11-74 down to 12-75.
$ws.Rows.Item(11).Insert()

# Populate the new row 11 with the inserted record's data.
$ws.Range("A11").Value2 = 1
$ws.Range("B11").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C11").Value2 = "Arica y Parinacota"
$ws.Range("D11").Value2 = 44831
$ws.Range("D11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E11").Value2 = 15
$ws.Range("F11").Value2 = 100112009
$ws.Range("G11").Value2 = "Acelga"
$ws.Range("H11").Value2 = "Sin especificar"
$ws.Range("I11").Value2 = "Primera"
$ws.Range("J11").Value2 = 250
$ws.Range("K11").Value2 = 1500
$ws.Range("L11").Value2 = 2000
$ws.Range("M11").Value2 = 1750
$ws.Range("N11").Value2 = "$/atado 2,5 a 3 kilos"
$ws.Range("O11").Value2 = "Región de Arica y Parinacota"
$ws.Range("P11").Value2 = 583
$ws.Range("Q11").Value2 = 3
$ws.Range("R11").Value2 = "Hortaliza"
